$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) YDS sheet: append this week's per-play yardage log to the four running
#    logs (OFF rush, OFF pass, DEF rush, DEF pass).
# ---------------------------------------------------------------------------
$ydsSheet = $wb.Worksheets.Item("YDS")

function Append-PlayLog($sheet, $cellRef, $newNumbers) {
    $cell = $sheet.Range($cellRef)
    $existing = $cell.Text
    $cell.Value = ($existing + " " + $newNumbers)
}

Append-PlayLog $ydsSheet "B2" "2 4 9 8 4 0 8 -2 0 13 0 11 4 9 3 -2 1 6 3 4 0 5 2 7 1 8 2"
Append-PlayLog $ydsSheet "B3" "7 7 5 10 24 6 5 7 14 22 4 10 9 5 40 8 9 25 16 13 1 18 11 12 8 13 8 61 9"
Append-PlayLog $ydsSheet "C2" "6 3 15 15 4 1 0 8 3 10 3 0 4 1 0 4 4 2 3"
Append-PlayLog $ydsSheet "C3" "14 48 8 2 21 -4 32 15 -2 19 5 5 6 12"

# ---------------------------------------------------------------------------
# 2) ST sheet: append this week's entries to the six running logs
#    (TB/D/RA/RM yardage + return logs).
# ---------------------------------------------------------------------------
$stSheet = $wb.Worksheets.Item("ST")

Append-PlayLog $stSheet "B4" "63"
Append-PlayLog $stSheet "B5" "0"
Append-PlayLog $stSheet "B6" "20"
Append-PlayLog $stSheet "D3" "47 54 62 51"
Append-PlayLog $stSheet "D4" "7 7 0 0"
Append-PlayLog $stSheet "D5" "21 0 0 0 0 0"

# ---------------------------------------------------------------------------
# 3) Season-total tables: overwrite the aggregate numbers on each sheet with
#    the updated totals (Week 16 logged + season simulated from Week 17).
# ---------------------------------------------------------------------------
$offSheet = $wb.Worksheets.Item("OFF")
$offSheet.Range("C2").Value = 203
$offSheet.Range("F2").Value = 58
$offSheet.Range("G2").Value = 51
$offSheet.Range("H2").Value = 4
$offSheet.Range("I2").Value = 6
$offSheet.Range("J2").Value = 30
$offSheet.Range("L2").Value = 241
$offSheet.Range("M2").Value = 164
$offSheet.Range("O2").Value = 19
$offSheet.Range("P2").Value = 9
$offSheet.Range("Q2").Value = 452
$offSheet.Range("C3").Value = 176
$offSheet.Range("D3").Value = 7
$offSheet.Range("E3").Value = 34
$offSheet.Range("F3").Value = 127
$offSheet.Range("G3").Value = 42
$offSheet.Range("H3").Value = 19
$offSheet.Range("I3").Value = 55
$offSheet.Range("J3").Value = 61
$offSheet.Range("N3").Value = 12

$defSheet = $wb.Worksheets.Item("DEF")
$defSheet.Range("C2").Value = 172
$defSheet.Range("E2").Value = 13
$defSheet.Range("F2").Value = 57
$defSheet.Range("J2").Value = 22
$defSheet.Range("L2").Value = 240
$defSheet.Range("M2").Value = 139
$defSheet.Range("O2").Value = 20
$defSheet.Range("P2").Value = 11
$defSheet.Range("Q2").Value = 432
$defSheet.Range("C3").Value = 179
$defSheet.Range("F3").Value = 104
$defSheet.Range("G3").Value = 33
$defSheet.Range("H3").Value = 28
$defSheet.Range("I3").Value = 60
$defSheet.Range("J3").Value = 51
$defSheet.Range("N3").Value = 21

$stSheet.Range("B2").Value = 93
$stSheet.Range("D2").Value = 58
$stSheet.Range("F2").Value = 350
$stSheet.Range("G2").Value = 336
$stSheet.Range("B3").Value = 78

$turnsSheet = $wb.Worksheets.Item("TURNS")
$turnsSheet.Range("D2").Value = 8
$turnsSheet.Range("E2").Value = 4
$turnsSheet.Range("D3").Value = 6
$turnsSheet.Range("E3").Value = 3

$penSheet = $wb.Worksheets.Item("PEN")
$penSheet.Range("D2").Value = 13
$penSheet.Range("D3").Value = 4
